# Regenerate merged AHB files:
#  - header row labels switch from the generic "_old"/"_new" suffixes
#    to the concrete version tags "_FV2410"/"_FV2504"
#  - the used range becomes a proper Excel Table ("Table1")
#  - the header row is frozen

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410",
    "Segment ID_FV2410", "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410", "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504", "Segmentgruppe_FV2504", "Segment_FV2504", "Datenelement_FV2504",
    "Segment ID_FV2504", "Code_FV2504", "Qualifier_FV2504", "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504", "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Convert A1:U58 into an Excel Table named "Table1"
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U58"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split above row 2, freeze the top pane)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
